$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends with a "r2_adj" results row (row 6: 0.65, 0.69,
# 0.24). We need to insert a new "r2" row above it with the same B/C
# values but a different D value (0.26), pushing "r2_adj" down to row 7.

# 1) Write the duplicated "r2_adj" row straight into new row 7 (this just
#    extends the used range by one row; row 6 is untouched so far).
$ws.Cells.Item(7,1).Value = "r2_adj"
$ws.Cells.Item(7,2).Value = 0.65
$ws.Cells.Item(7,3).Value = 0.6899999999999999
$ws.Cells.Item(7,4).Value = 0.24

# 2) Copy row 6's label formatting (bold, centered, bordered) onto the
#    new A7 label cell so it matches the other row labels in column A.
$ws.Cells.Item(6,1).Copy()
$ws.Cells.Item(7,1).PasteSpecial(-4122)

# 3) Turn the original row 6 into the new "r2" row: same label style,
#    same B/C values, but a new D value.
$ws.Cells.Item(6,1).Value = "r2"
$ws.Cells.Item(6,4).Value = 0.26
